$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above the current row 88 (old rows 88-99 shift down to 91-102)
$ws.Range("A88:T90").EntireRow.Insert()

# --- New row 88: Chirimoya, Especial, $/bandeja 10 kilos, Provincia de Limarí ---
$ws.Cells.Item(88,1).Value = 2
$ws.Cells.Item(88,2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(88,3).Value = "Coquimbo"
$ws.Cells.Item(88,4).Value = 44826
$ws.Cells.Item(88,5).Value = 4
$ws.Cells.Item(88,6).Value = "Fruta"
$ws.Cells.Item(88,7).Value = 100107
$ws.Cells.Item(88,8).Value = "Otros"
$ws.Cells.Item(88,9).Value = 100107002
$ws.Cells.Item(88,10).Value = "Chirimoya"
$ws.Cells.Item(88,11).Value = "Cultivar IV Región"
$ws.Cells.Item(88,12).Value = "Especial"
$ws.Cells.Item(88,13).Value = 400
$ws.Cells.Item(88,14).Value = 22000
$ws.Cells.Item(88,15).Value = 23000
$ws.Cells.Item(88,16).Value = 22500
$ws.Cells.Item(88,17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(88,18).Value = "Provincia de Limarí"
$ws.Cells.Item(88,19).Value = 2250
$ws.Cells.Item(88,20).Value = 10

# --- New row 89: Chirimoya, Primera, $/bandeja 10 kilos, Provincia de Limarí ---
$ws.Cells.Item(89,1).Value = 2
$ws.Cells.Item(89,2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(89,3).Value = "Coquimbo"
$ws.Cells.Item(89,4).Value = 44826
$ws.Cells.Item(89,5).Value = 4
$ws.Cells.Item(89,6).Value = "Fruta"
$ws.Cells.Item(89,7).Value = 100107
$ws.Cells.Item(89,8).Value = "Otros"
$ws.Cells.Item(89,9).Value = 100107002
$ws.Cells.Item(89,10).Value = "Chirimoya"
$ws.Cells.Item(89,11).Value = "Cultivar IV Región"
$ws.Cells.Item(89,12).Value = "Primera"
$ws.Cells.Item(89,13).Value = 300
$ws.Cells.Item(89,14).Value = 18000
$ws.Cells.Item(89,15).Value = 19000
$ws.Cells.Item(89,16).Value = 18500
$ws.Cells.Item(89,17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(89,18).Value = "Provincia de Limarí"
$ws.Cells.Item(89,19).Value = 1850
$ws.Cells.Item(89,20).Value = 10

# --- New row 90: Chirimoya, Segunda, $/bandeja 10 kilos, Provincia de Limarí ---
$ws.Cells.Item(90,1).Value = 2
$ws.Cells.Item(90,2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(90,3).Value = "Coquimbo"
$ws.Cells.Item(90,4).Value = 44826
$ws.Cells.Item(90,5).Value = 4
$ws.Cells.Item(90,6).Value = "Fruta"
$ws.Cells.Item(90,7).Value = 100107
$ws.Cells.Item(90,8).Value = "Otros"
$ws.Cells.Item(90,9).Value = 100107002
$ws.Cells.Item(90,10).Value = "Chirimoya"
$ws.Cells.Item(90,11).Value = "Cultivar IV Región"
$ws.Cells.Item(90,12).Value = "Segunda"
$ws.Cells.Item(90,13).Value = 240
$ws.Cells.Item(90,14).Value = 15000
$ws.Cells.Item(90,15).Value = 16000
$ws.Cells.Item(90,16).Value = 15500
$ws.Cells.Item(90,17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(90,18).Value = "Provincia de Limarí"
$ws.Cells.Item(90,19).Value = 1550
$ws.Cells.Item(90,20).Value = 10
